$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 48 (shifts existing rows 48:110 down to 49:111)
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly price-report record
$ws.Range("A48").Value = 7
$ws.Range("B48").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C48").Value = "Ñuble"
$ws.Range("D48").Value = 44893
$ws.Range("E48").Value = 16
$ws.Range("F48").Value = 100112021
$ws.Range("G48").Value = "Ají"
$ws.Range("H48").Value = "Americana (o)"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 60
$ws.Range("K48").Value = 16000
$ws.Range("L48").Value = 17000
$ws.Range("M48").Value = 16500
$ws.Range("N48").Value = "`$/caja 15 kilos"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 1100
$ws.Range("Q48").Value = 15
$ws.Range("R48").Value = "Hortaliza"
